$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 5, shifting rows 5-20 down to 6-21
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row with the new entry
$ws.Cells.Item(5, 1).Value = "m4/~`$Boss.xlsx"
$ws.Cells.Item(5, 2).Value = "~`$Boss"
